# Insert a new daily data row before row 901 (shifts rows 901:942 down to 902:943)
# and populate it with the new auto-pushed data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(901).Insert()

$ws.Range("A901").Value = "2026/02/28"
$ws.Range("B901").Value = "土"
$ws.Range("C901").Value = 7
$ws.Range("D901").Value = 201
